$d = $word.ActiveDocument

# 1. Change the text of the second run in the first paragraph from " - V1" to " "
$d.Content.Find.Execute(" – V1", $false, $false, $false, $false, $false,
                         $true, 1, $false, " ", 2)

# 2. Find the "Requisitos do sistema:" paragraph and color it red (both paragraph
#    mark and the run text)
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

# Color the run text red
$p2.Range.Font.Color = 255  # wdColorRed (0x0000FF in BGR = 255)

# Color the paragraph mark (pPr/rPr) red as well: include the paragraph mark
# in the range so its rPr gets the color property.
$rng = $p2.Range
$rng.MoveEnd(1, 1)  # extend to include the paragraph mark
$rng.Font.Color = 255

# 3. Move the _GoBack bookmark from end of paragraph 1 to end of paragraph 2
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$p2end = $p2.Range
$p2end.Collapse(0)  # collapse to end (wdCollapseEnd = 0)
$d.Bookmarks.Add("_GoBack", $p2end)
